$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, pushing existing rows 17-76 down to 18-77.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new record's data.
$ws.Cells.Item(17, 1).Value = 11
$ws.Cells.Item(17, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(17, 3).Value = "Bíobío"
$ws.Cells.Item(17, 4).Value = 44614
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 8
$ws.Cells.Item(17, 6).Value = 100112021
$ws.Cells.Item(17, 7).Value = "Ají"
$ws.Cells.Item(17, 8).Value = "Chilena(o)"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 24000
$ws.Cells.Item(17, 12).Value = 25000
$ws.Cells.Item(17, 13).Value = 24500
$ws.Cells.Item(17, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 980
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"
